$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header text in F1: "Draft" -> "Drafting of manuscript"
$ws.Range("F1").Value = "Drafting of manuscript"

# Update the active cell selection from E2 to F2
$ws.Range("F2").Select()
